# Update "想去人数" (interested count) values on both the "展览" and
# "全部类型" sheets to reflect newly generated output.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F11").Value = 76
$ws1.Range("F14").Value = 562

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 76
$ws4.Range("F15").Value = 562
